$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the Sub Category label so the KPI exclusion lookup matches the
# expected key name (the template's exclusion logic keys off this exact
# string).
$ws.Range("B2").Value = "sub_category"

# Nudge the explicit column widths back toward the values the workbook
# carried after the fix was re-saved (best match achievable through the
# exposed ColumnWidth API, which only resolves to whole pixels).
$ws.Cells.Item(1, 1).ColumnWidth = 18.450067476383268
$ws.Cells.Item(1, 2).ColumnWidth = 14.591767881241566
$ws.Cells.Item(1, 3).ColumnWidth = 207.83468286099867
$ws.Cells.Item(1, 4).ColumnWidth = 13.948043184885266
$ws.Cells.Item(1, 5).ColumnWidth = 22.518893387314467

# Restore cursor position as left by the author when the fix was saved.
$ws.Range("C24").Select()
